$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 20.00768154687529
$ws.Cells.Item(2, 3).Value = 12.58822934047073
$ws.Cells.Item(2, 4).Value = 15.11417044970839
$ws.Cells.Item(2, 5).Value = 16.54545837451031
$ws.Cells.Item(2, 7).Value = 58.99718260621379
$ws.Cells.Item(2, 8).Value = 21.1626555222491
$ws.Cells.Item(2, 10).Value = 9.459625741028404

$ws.Cells.Item(3, 2).Value = 19.58021695603266
$ws.Cells.Item(3, 3).Value = 12.19389734686851
$ws.Cells.Item(3, 4).Value = 15.06430076538702
$ws.Cells.Item(3, 5).Value = 16.49689304671548
$ws.Cells.Item(3, 7).Value = 58.31384738202864
$ws.Cells.Item(3, 8).Value = 21.11384465096936
$ws.Cells.Item(3, 10).Value = 9.473865552789846

$ws.Cells.Item(4, 2).Value = 19.31975601849749
$ws.Cells.Item(4, 3).Value = 11.94966845093827
$ws.Cells.Item(4, 4).Value = 15.0374347015375
$ws.Cells.Item(4, 5).Value = 16.47120512415311
$ws.Cells.Item(4, 7).Value = 57.90922598828651
$ws.Cells.Item(4, 8).Value = 21.08887383393095
$ws.Cells.Item(4, 10).Value = 9.484223070089648

$ws.Cells.Item(5, 2).Value = 19.21428940576688
$ws.Cells.Item(5, 3).Value = 11.84979418922796
$ws.Cells.Item(5, 4).Value = 15.02743548200475
$ws.Cells.Item(5, 5).Value = 16.46178006826518
$ws.Cells.Item(5, 7).Value = 57.74825038128382
$ws.Cells.Item(5, 8).Value = 21.07995436010414
$ws.Cells.Item(5, 10).Value = 9.488849062884592

$ws.Cells.Item(6, 2).Value = 19.19682241499803
$ws.Cells.Item(6, 3).Value = 11.8331945695995
$ws.Cells.Item(6, 4).Value = 15.02583255326654
$ws.Cells.Item(6, 5).Value = 16.46027813577395
$ws.Cells.Item(6, 7).Value = 57.72176094460323
$ws.Cells.Item(6, 8).Value = 21.07854912061721
$ws.Cells.Item(6, 10).Value = 9.489641658573525

$ws.Cells.Item(7, 2).Value = 19.31833070707778
$ws.Cells.Item(7, 3).Value = 11.94832267752639
$ws.Cells.Item(7, 4).Value = 15.03729600073778
$ws.Cells.Item(7, 5).Value = 16.47107378696281
$ws.Cells.Item(7, 7).Value = 57.90703898688466
$ws.Cells.Item(7, 8).Value = 21.08874845783949
$ws.Cells.Item(7, 10).Value = 9.48428381806086

$ws.Cells.Item(8, 2).Value = 19.85997663872913
$ws.Cells.Item(8, 3).Value = 12.45280608937197
$ws.Cells.Item(8, 4).Value = 15.09619949322178
$ws.Cells.Item(8, 5).Value = 16.52785838832906
$ws.Cells.Item(8, 7).Value = 58.75856004702167
$ws.Cells.Item(8, 8).Value = 21.14478784381458
$ws.Cells.Item(8, 10).Value = 9.464200199401429

$ws.Cells.Item(9, 2).Value = 20.9307648531812
$ws.Cells.Item(9, 3).Value = 13.41802590897267
$ws.Cells.Item(9, 4).Value = 15.24123172417699
$ws.Cells.Item(9, 5).Value = 16.67176290064277
$ws.Cells.Item(9, 7).Value = 60.53987609377348
$ws.Cells.Item(9, 8).Value = 21.29432979461035
$ws.Cells.Item(9, 10).Value = 9.437653823622117

$ws.Cells.Item(10, 2).Value = 21.71333330338366
$ws.Cells.Item(10, 3).Value = 14.10355742886801
$ws.Cells.Item(10, 4).Value = 15.36536168320852
$ws.Cells.Item(10, 5).Value = 16.79695320806165
$ws.Cells.Item(10, 7).Value = 61.90649848972031
$ws.Cells.Item(10, 8).Value = 21.4282708774779
$ws.Cells.Item(10, 10).Value = 9.426014748029729

$ws.Cells.Item(11, 2).Value = 22.06658891080751
$ws.Cells.Item(11, 3).Value = 14.40866894868194
$ws.Cells.Item(11, 4).Value = 15.42553593177851
$ws.Cells.Item(11, 5).Value = 16.85802802746422
$ws.Cells.Item(11, 7).Value = 62.53852720154283
$ws.Cells.Item(11, 8).Value = 21.49438414607384
$ws.Cells.Item(11, 10).Value = 9.422434599479889

$ws.Cells.Item(12, 2).Value = 22.19982017554962
$ws.Cells.Item(12, 3).Value = 14.52311969279547
$ws.Cells.Item(12, 4).Value = 15.44884434555595
$ws.Cells.Item(12, 5).Value = 16.88173815097268
$ws.Cells.Item(12, 7).Value = 62.77915446448527
$ws.Cells.Item(12, 8).Value = 21.52015862718561
$ws.Cells.Item(12, 10).Value = 9.421325877020418

$ws.Cells.Item(13, 2).Value = 22.17115262437832
$ws.Cells.Item(13, 3).Value = 14.49852084611271
$ws.Cells.Item(13, 4).Value = 15.443801447386
$ws.Cells.Item(13, 5).Value = 16.8766060306503
$ws.Cells.Item(13, 7).Value = 62.72727675418056
$ws.Cells.Item(13, 8).Value = 21.5145748880008
$ws.Cells.Item(13, 10).Value = 9.421553668252274

$ws.Cells.Item(14, 2).Value = 22.0775615281028
$ws.Cells.Item(14, 3).Value = 14.41810733045351
$ws.Cells.Item(14, 4).Value = 15.42744314080653
$ws.Cells.Item(14, 5).Value = 16.8599670663578
$ws.Cells.Item(14, 7).Value = 62.55829898786027
$ws.Cells.Item(14, 8).Value = 21.49648986242155
$ws.Cells.Item(14, 10).Value = 9.422338431568477

$ws.Cells.Item(15, 2).Value = 22.02015993603804
$ws.Cells.Item(15, 3).Value = 14.36870667206183
$ws.Cells.Item(15, 4).Value = 15.41749079899456
$ws.Cells.Item(15, 5).Value = 16.84985073122454
$ws.Cells.Item(15, 7).Value = 62.45495761930009
$ws.Cells.Item(15, 8).Value = 21.48550826577697
$ws.Cells.Item(15, 10).Value = 9.422851301036207

$ws.Cells.Item(16, 2).Value = 21.69017949706965
$ws.Cells.Item(16, 3).Value = 14.083471570304
$ws.Cells.Item(16, 4).Value = 15.36150280083351
$ws.Cells.Item(16, 5).Value = 16.79304393500976
$ws.Cells.Item(16, 7).Value = 61.86538537138514
$ws.Cells.Item(16, 8).Value = 21.42405398697107
$ws.Cells.Item(16, 10).Value = 9.426283260327754

$ws.Cells.Item(17, 2).Value = 21.48693964210493
$ws.Cells.Item(17, 3).Value = 13.90667394357934
$ws.Cells.Item(17, 4).Value = 15.32809765096748
$ws.Cells.Item(17, 5).Value = 16.75924400407575
$ws.Cells.Item(17, 7).Value = 61.50621249752577
$ws.Cells.Item(17, 8).Value = 21.38767701247813
$ws.Cells.Item(17, 10).Value = 9.428828133557261

$ws.Cells.Item(18, 2).Value = 21.36979180344068
$ws.Cells.Item(18, 3).Value = 13.80435650377699
$ws.Cells.Item(18, 4).Value = 15.30923354262927
$ws.Cells.Item(18, 5).Value = 16.74019216423924
$ws.Cells.Item(18, 7).Value = 61.30061134808093
$ws.Cells.Item(18, 8).Value = 21.36724223467791
$ws.Cells.Item(18, 10).Value = 9.430453233806277

$ws.Cells.Item(19, 2).Value = 21.33008924935979
$ws.Cells.Item(19, 3).Value = 13.76960953433759
$ws.Cells.Item(19, 4).Value = 15.30290686195852
$ws.Cells.Item(19, 5).Value = 16.73380865913581
$ws.Cells.Item(19, 7).Value = 61.23117356183446
$ws.Cells.Item(19, 8).Value = 21.36040741180098
$ws.Cells.Item(19, 10).Value = 9.431031162912811

$ws.Cells.Item(20, 2).Value = 21.50860178934475
$ws.Cells.Item(20, 3).Value = 13.92556026403072
$ws.Cells.Item(20, 4).Value = 15.3316175777552
$ws.Cells.Item(20, 5).Value = 16.76280188066388
$ws.Cells.Item(20, 7).Value = 61.54434647542265
$ws.Cells.Item(20, 8).Value = 21.39149890965883
$ws.Cells.Item(20, 10).Value = 9.428540524157016

$ws.Cells.Item(21, 2).Value = 22.10506724109756
$ws.Cells.Item(21, 3).Value = 14.44175715431238
$ws.Cells.Item(21, 4).Value = 15.43223390706966
$ws.Cells.Item(21, 5).Value = 16.8648386177011
$ws.Cells.Item(21, 7).Value = 62.6078983616537
$ws.Cells.Item(21, 8).Value = 21.50178187894392
$ws.Cells.Item(21, 10).Value = 9.422101220713088

$ws.Cells.Item(22, 2).Value = 22.49168951042718
$ws.Cells.Item(22, 3).Value = 14.77273131417206
$ws.Cells.Item(22, 4).Value = 15.50102766424759
$ws.Cells.Item(22, 5).Value = 16.93491460038646
$ws.Cells.Item(22, 7).Value = 63.31042910498945
$ws.Cells.Item(22, 8).Value = 21.57816074098246
$ws.Cells.Item(22, 10).Value = 9.419332642199665

$ws.Cells.Item(23, 2).Value = 22.28568054661751
$ws.Cells.Item(23, 3).Value = 14.59670503022621
$ws.Cells.Item(23, 4).Value = 15.46403739564346
$ws.Cells.Item(23, 5).Value = 16.89720744083482
$ws.Cells.Item(23, 7).Value = 62.93485944103413
$ws.Cells.Item(23, 8).Value = 21.5370046667881
$ws.Cells.Item(23, 10).Value = 9.420678403450342

$ws.Cells.Item(24, 2).Value = 21.49880926380713
$ws.Cells.Item(24, 3).Value = 13.91702384914149
$ws.Cells.Item(24, 4).Value = 15.33002515600178
$ws.Cells.Item(24, 5).Value = 16.76119217964427
$ws.Cells.Item(24, 7).Value = 61.52710331040736
$ws.Cells.Item(24, 8).Value = 21.38976953802194
$ws.Cells.Item(24, 10).Value = 9.428670047749819

$ws.Cells.Item(25, 2).Value = 20.64115946994467
$ws.Cells.Item(25, 3).Value = 13.16045945411401
$ws.Cells.Item(25, 4).Value = 15.19887547553332
$ws.Cells.Item(25, 5).Value = 16.62938293244216
$ws.Cells.Item(25, 7).Value = 60.04706351692766
$ws.Cells.Item(25, 8).Value = 21.24963511271217
$ws.Cells.Item(25, 10).Value = 9.443456752685787
